$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 36.94436433333333
$ws.Range("H2").Value = 110.833093
$ws.Range("I2").Value = 0.8328964975864823
$ws.Range("J2").Value = 0.8328964975864824
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 19.11595033333333
$ws.Range("N2").Value = 57.347851
$ws.Range("O2").Value = 0.6851940154453416
$ws.Range("P2").Value = 0.6851940154453418
$ws.Range("Q2").Value = 706.2266336925715
$ws.Range("R2").Value = 6356.039703233143
$ws.Range("S2").Value = 0.5706956956316431
$ws.Range("T2").Value = 0.5706956956316432

$ws.Range("G3").Value = 36.94436433333333
$ws.Range("H3").Value = 110.833093
$ws.Range("I3").Value = 0.8328964975864823
$ws.Range("J3").Value = 0.8328964975864824
$ws.Range("M3").Value = 4.865208333333334
$ws.Range("O3").Value = 0.1743890089566637
$ws.Range("P3").Value = 0.1743890089566637
$ws.Range("Q3").Value = 179.7420292242361
$ws.Range("R3").Value = 1617.678263018125
$ws.Range("S3").Value = 0.1452479947775829
$ws.Range("T3").Value = 0.1452479947775829

$ws.Range("G4").Value = 36.94436433333333
$ws.Range("H4").Value = 110.833093
$ws.Range("I4").Value = 0.8328964975864823
$ws.Range("J4").Value = 0.8328964975864824
$ws.Range("M4").Value = 3.917436333333333
$ws.Range("N4").Value = 11.752309
$ws.Range("O4").Value = 0.1404169755979945
$ws.Range("P4").Value = 0.1404169755979946
$ws.Range("Q4").Value = 144.7271951513041
$ws.Range("R4").Value = 1302.544756361737
$ws.Range("S4").Value = 0.1169528071772562
$ws.Range("T4").Value = 0.1169528071772562

$ws.Range("I5").Value = 0.07608399754092349
$ws.Range("J5").Value = 0.07608399754092349
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 19.11595033333333
$ws.Range("N5").Value = 57.347851
$ws.Range("O5").Value = 0.6851940154453416
$ws.Range("P5").Value = 0.6851940154453418
$ws.Range("Q5").Value = 64.51287238798967
$ws.Range("R5").Value = 580.6158514919069
$ws.Range("S5").Value = 0.05213229978619886
$ws.Range("T5").Value = 0.05213229978619887

$ws.Range("I6").Value = 0.07608399754092349
$ws.Range("J6").Value = 0.07608399754092349
$ws.Range("M6").Value = 4.865208333333334
$ws.Range("O6").Value = 0.1743890089566637
$ws.Range("P6").Value = 0.1743890089566637
$ws.Range("S6").Value = 0.01326821292862288
$ws.Range("T6").Value = 0.01326821292862289

$ws.Range("I7").Value = 0.07608399754092349
$ws.Range("J7").Value = 0.07608399754092349
$ws.Range("M7").Value = 3.917436333333333
$ws.Range("N7").Value = 11.752309
$ws.Range("O7").Value = 0.1404169755979945
$ws.Range("P7").Value = 0.1404169755979946
$ws.Range("Q7").Value = 13.22063856902367
$ws.Range("R7").Value = 118.985747121213
$ws.Range("S7").Value = 0.01068348482610173
$ws.Range("T7").Value = 0.01068348482610173

$ws.Range("G8").Value = 4.037305666666668
$ws.Range("H8").Value = 12.111917
$ws.Range("I8").Value = 0.09101950487259411
$ws.Range("J8").Value = 0.09101950487259411
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 19.11595033333333
$ws.Range("N8").Value = 57.347851
$ws.Range("O8").Value = 0.6851940154453416
$ws.Range("P8").Value = 0.6851940154453418
$ws.Range("Q8").Value = 77.17693460448524
$ws.Range("R8").Value = 694.5924114403671
$ws.Range("S8").Value = 0.06236602002749959
$ws.Range("T8").Value = 0.06236602002749961

$ws.Range("G9").Value = 4.037305666666668
$ws.Range("H9").Value = 12.111917
$ws.Range("I9").Value = 0.09101950487259411
$ws.Range("J9").Value = 0.09101950487259411
$ws.Range("M9").Value = 4.865208333333334
$ws.Range("O9").Value = 0.1743890089566637
$ws.Range("P9").Value = 0.1743890089566637
$ws.Range("Q9").Value = 19.64233317368056
$ws.Range("R9").Value = 176.7809985631251
$ws.Range("S9").Value = 0.01587280125045791
$ws.Range("T9").Value = 0.01587280125045791

$ws.Range("G10").Value = 4.037305666666668
$ws.Range("H10").Value = 12.111917
$ws.Range("I10").Value = 0.09101950487259411
$ws.Range("J10").Value = 0.09101950487259411
$ws.Range("M10").Value = 3.917436333333333
$ws.Range("N10").Value = 11.752309
$ws.Range("O10").Value = 0.1404169755979945
$ws.Range("P10").Value = 0.1404169755979946
$ws.Range("Q10").Value = 15.81588790737256
$ws.Range("R10").Value = 142.342991166353
$ws.Range("S10").Value = 0.01278068359463659
$ws.Range("T10").Value = 0.01278068359463659
